$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace header titles with Persian translations
$ws.Range("A1").Value = "ستون اول"
$ws.Range("B1").Value = "ستون دوم"
$ws.Range("C1").Value = "ستون سوم"

# Update the active selection to C2
$ws.Range("C2").Select()
